$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44315
$ws.Range("I2").Value = 'Especial'
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 30000
$ws.Range("M2").Value = 30000
$ws.Range("N2").Value = '$/caja 20 kilos empedrada'
$ws.Range("P2").Value = 1500
$ws.Range("Q2").Value = 20

# Row 3
$ws.Range("D3").Value = 44315
$ws.Range("N3").Value = '$/caja 15 kilos granel'

# Row 4
$ws.Range("D4").Value = 44293
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 25000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 25000
$ws.Range("N4").Value = '$/caja 15 kilos empedrada'
$ws.Range("P4").Value = 1667
$ws.Range("Q4").Value = 15

# Row 5
$ws.Range("D5").Value = 44313
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("P5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44313
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 30000
$ws.Range("N6").Value = '$/caja 20 kilos empedrada'
$ws.Range("P6").Value = 1500
$ws.Range("Q6").Value = 20

# Row 7
$ws.Range("D7").Value = 44280
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 25000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 25000
$ws.Range("N7").Value = '$/caja 18 kilos empedrada'
$ws.Range("P7").Value = 1389
$ws.Range("Q7").Value = 18

# Row 8
$ws.Range("D8").Value = 44285
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 25000
$ws.Range("N8").Value = '$/caja 18 kilos empedrada'
$ws.Range("P8").Value = 1389
$ws.Range("Q8").Value = 18
